$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Efna1"
$ws.Range("C2").Value = "Epha3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 14.02618833333333
$ws.Range("H2").Value = 42.078565
$ws.Range("I2").Value = 0.806325281849088
$ws.Range("J2").Value = 0.8172785134657441
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.107177
$ws.Range("N2").Value = 0.321531
$ws.Range("O2").Value = 0.003526763356587491
$ws.Range("P2").Value = 0.003549676734010809
$ws.Range("Q2").Value = 1.503284787001667
$ws.Range("R2").Value = 13.529563083015
$ws.Range("S2").Value = 0.002843718457515444
$ws.Range("T2").Value = 0.002901074524456291

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Efna1"
$ws.Range("C3").Value = "Epha3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 14.02618833333333
$ws.Range("H3").Value = 42.078565
$ws.Range("I3").Value = 0.806325281849088
$ws.Range("J3").Value = 0.8172785134657441
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 29.65321533333333
$ws.Range("N3").Value = 88.95964599999999
$ws.Range("O3").Value = 0.9757678722356318
$ws.Range("P3").Value = 0.9821074349659524
$ws.Range("Q3").Value = 415.921582954221
$ws.Range("R3").Value = 3743.29424658799
$ws.Range("S3").Value = 0.7867863045996807
$ws.Range("T3").Value = 0.8026553045126286

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Efna1"
$ws.Range("C4").Value = "Epha3"
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 14.02618833333333
$ws.Range("H4").Value = 42.078565
$ws.Range("I4").Value = 0.806325281849088
$ws.Range("J4").Value = 0.8172785134657441
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.04072766666666667
$ws.Range("N4").Value = 0.122183
$ws.Range("O4").Value = 0.001340183457265176
$ws.Range("P4").Value = 0.001348890627627329
$ws.Range("Q4").Value = 0.5712539230438889
$ws.Range("R4").Value = 5.141285307395
$ws.Range("S4").Value = 0.001080623803908828
$ws.Range("T4").Value = 0.001102419326975138

$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Efna1"
$ws.Range("C5").Value = "Epha3"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 14.02618833333333
$ws.Range("H5").Value = 42.078565
$ws.Range("I5").Value = 0.806325281849088
$ws.Range("J5").Value = 0.8172785134657441
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.5885005
$ws.Range("N5").Value = 1.177001
$ws.Range("O5").Value = 0.01936518095051565
$ws.Range("P5").Value = 0.01299399767240936
$ws.Range("Q5").Value = 8.254418847260832
$ws.Range("R5").Value = 49.52651308356499
$ws.Range("S5").Value = 0.01561463498798313
$ws.Range("T5").Value = 0.01061971510168406

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Efna1"
$ws.Range("C6").Value = "Epha3"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 2.483247333333333
$ws.Range("H6").Value = 7.449742
$ws.Range("I6").Value = 0.1427547569137158
$ws.Range("J6").Value = 0.1446939568272663
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.107177
$ws.Range("N6").Value = 0.321531
$ws.Range("O6").Value = 0.003526763356587491
$ws.Range("P6").Value = 0.003549676734010809
$ws.Range("Q6").Value = 0.2661469994446667
$ws.Range("R6").Value = 2.395322995002
$ws.Range("S6").Value = 0.0005034622456618476
$ws.Range("T6").Value = 0.0005136167721017117

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Efna1"
$ws.Range("C7").Value = "Epha3"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 2.483247333333333
$ws.Range("H7").Value = 7.449742
$ws.Range("I7").Value = 0.1427547569137158
$ws.Range("J7").Value = 0.1446939568272663
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 29.65321533333333
$ws.Range("N7").Value = 88.95964599999999
$ws.Range("O7").Value = 0.9757678722356318
$ws.Range("P7").Value = 0.9821074349659524
$ws.Range("Q7").Value = 73.6362679012591
$ws.Range("R7").Value = 662.7264111113319
$ws.Range("S7").Value = 0.1392955054052113
$ws.Range("T7").Value = 0.1421050107947008

$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Efna1"
$ws.Range("C8").Value = "Epha3"
$ws.Range("D8").Value = "Inflammatory-Mac"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 2.483247333333333
$ws.Range("H8").Value = 7.449742
$ws.Range("I8").Value = 0.1427547569137158
$ws.Range("J8").Value = 0.1446939568272663
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.04072766666666667
$ws.Range("N8").Value = 0.122183
$ws.Range("O8").Value = 0.001340183457265176
$ws.Range("P8").Value = 0.001348890627627329
$ws.Range("Q8").Value = 0.1011368696428889
$ws.Range("R8").Value = 0.910231826786
$ws.Range("S8").Value = 0.0001913175636616734
$ws.Range("T8").Value = 0.0001951763222386129

$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Efna1"
$ws.Range("C9").Value = "Epha3"
$ws.Range("D9").Value = "MuSCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 2.483247333333333
$ws.Range("H9").Value = 7.449742
$ws.Range("I9").Value = 0.1427547569137158
$ws.Range("J9").Value = 0.1446939568272663
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.5885005
$ws.Range("N9").Value = 1.177001
$ws.Range("O9").Value = 0.01936518095051565
$ws.Range("P9").Value = 0.01299399767240936
$ws.Range("Q9").Value = 1.461392297290333
$ws.Range("R9").Value = 8.768353783741999
$ws.Range("S9").Value = 0.002764471699180982
$ws.Range("T9").Value = 0.0018801529382252

$ws.Range("A10").Value = "Inflammatory-Mac"
$ws.Range("B10").Value = "Efna1"
$ws.Range("C10").Value = "Epha3"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.1863673333333333
$ws.Range("H10").Value = 0.559102
$ws.Range("I10").Value = 0.01071372271683668
$ws.Range("J10").Value = 0.01085925937435662
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.107177
$ws.Range("N10").Value = 0.321531
$ws.Range("O10").Value = 0.003526763356587491
$ws.Range("P10").Value = 0.003549676734010809
$ws.Range("Q10").Value = 0.01997429168466667
$ws.Range("R10").Value = 0.179768625162
$ws.Range("S10").Value = 0.00003778476469037858
$ws.Range("T10").Value = 0.00003854686034974248

$ws.Range("A11").Value = "Inflammatory-Mac"
$ws.Range("B11").Value = "Efna1"
$ws.Range("C11").Value = "Epha3"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.1863673333333333
$ws.Range("H11").Value = 0.559102
$ws.Range("I11").Value = 0.01071372271683668
$ws.Range("J11").Value = 0.01085925937435662
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 29.65321533333333
$ws.Range("N11").Value = 88.95964599999999
$ws.Range("O11").Value = 0.9757678722356318
$ws.Range("P11").Value = 0.9821074349659524
$ws.Range("Q11").Value = 5.526390666432444
$ws.Range("R11").Value = 49.737515997892
$ws.Range("S11").Value = 0.01045410641913028
$ws.Range("T11").Value = 0.01066495936977936

$ws.Range("A12").Value = "Inflammatory-Mac"
$ws.Range("B12").Value = "Efna1"
$ws.Range("C12").Value = "Epha3"
$ws.Range("D12").Value = "Inflammatory-Mac"
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.1863673333333333
$ws.Range("H12").Value = 0.559102
$ws.Range("I12").Value = 0.01071372271683668
$ws.Range("J12").Value = 0.01085925937435662
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.04072766666666667
$ws.Range("N12").Value = 0.122183
$ws.Range("O12").Value = 0.001340183457265176
$ws.Range("P12").Value = 0.001348890627627329
$ws.Range("Q12").Value = 0.007590306629555556
$ws.Range("R12").Value = 0.06831275966599999
$ws.Range("S12").Value = 0.00001435835395083064
$ws.Range("T12").Value = 0.00001464795319304386

$ws.Range("A13").Value = "Inflammatory-Mac"
$ws.Range("B13").Value = "Efna1"
$ws.Range("C13").Value = "Epha3"
$ws.Range("D13").Value = "MuSCs"
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.1863673333333333
$ws.Range("H13").Value = 0.559102
$ws.Range("I13").Value = 0.01071372271683668
$ws.Range("J13").Value = 0.01085925937435662
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.5885005
$ws.Range("N13").Value = 1.177001
$ws.Range("O13").Value = 0.01936518095051565
$ws.Range("P13").Value = 0.01299399767240936
$ws.Range("Q13").Value = 0.1096772688503333
$ws.Range("R13").Value = 0.658063613102
$ws.Range("S13").Value = 0.0002074731790651925
$ws.Range("T13").Value = 0.0001411051910344795

$ws.Range("A14").Value = "MuSCs"
$ws.Range("B14").Value = "Efna1"
$ws.Range("C14").Value = "Epha3"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.6993955000000001
$ws.Range("H14").Value = 1.398791
$ws.Range("I14").Value = 0.04020623852035952
$ws.Range("J14").Value = 0.02716827033263282
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 0.107177
$ws.Range("N14").Value = 0.321531
$ws.Range("O14").Value = 0.003526763356587491
$ws.Range("P14").Value = 0.003549676734010809
$ws.Range("Q14").Value = 0.07495911150350001
$ws.Range("R14").Value = 0.449754669021
$ws.Range("S14").Value = 0.0001417978887198204
$ws.Range("T14").Value = 0.00009643857710306282

$ws.Range("A15").Value = "MuSCs"
$ws.Range("B15").Value = "Efna1"
$ws.Range("C15").Value = "Epha3"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.6993955000000001
$ws.Range("H15").Value = 1.398791
$ws.Range("I15").Value = 0.04020623852035952
$ws.Range("J15").Value = 0.02716827033263282
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 29.65321533333333
$ws.Range("N15").Value = 88.95964599999999
$ws.Range("O15").Value = 0.9757678722356318
$ws.Range("P15").Value = 0.9821074349659524
$ws.Range("Q15").Value = 20.73932536466434
$ws.Range("R15").Value = 124.435952187986
$ws.Range("S15").Value = 0.0392319558116095
$ws.Range("T15").Value = 0.02668216028884361

$ws.Range("A16").Value = "MuSCs"
$ws.Range("B16").Value = "Efna1"
$ws.Range("C16").Value = "Epha3"
$ws.Range("D16").Value = "Inflammatory-Mac"
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.6993955000000001
$ws.Range("H16").Value = 1.398791
$ws.Range("I16").Value = 0.04020623852035952
$ws.Range("J16").Value = 0.02716827033263282
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.04072766666666667
$ws.Range("N16").Value = 0.122183
$ws.Range("O16").Value = 0.001340183457265176
$ws.Range("P16").Value = 0.001348890627627329
$ws.Range("Q16").Value = 0.02848474679216667
$ws.Range("R16").Value = 0.170908480753
$ws.Range("S16").Value = 0.00005388373574384372
$ws.Range("T16").Value = 0.00003664702522053403

$ws.Range("A17").Value = "MuSCs"
$ws.Range("B17").Value = "Efna1"
$ws.Range("C17").Value = "Epha3"
$ws.Range("D17").Value = "MuSCs"
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.6993955000000001
$ws.Range("H17").Value = 1.398791
$ws.Range("I17").Value = 0.04020623852035952
$ws.Range("J17").Value = 0.02716827033263282
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.5885005
$ws.Range("N17").Value = 1.177001
$ws.Range("O17").Value = 0.01936518095051565
$ws.Range("P17").Value = 0.01299399767240936
$ws.Range("Q17").Value = 0.41159460144775
$ws.Range("R17").Value = 1.646378405791
$ws.Range("S17").Value = 0.0007786010842863549
$ws.Range("T17").Value = 0.0003530244414656193

Write-Output "done"
